# Simplify empty "Compact"-styled table-cell paragraphs down to bare,
# styleless empty paragraphs (<w:p/>), matching the target edit: cells
# in the Schedule table that contain no visible text (e.g. unfilled
# "Due" column entries) should no longer carry an explicit
# pPr/pStyle="Compact" - they become plain empty paragraphs.

$d = $word.ActiveDocument

foreach ($t in $d.Tables) {
    for ($r = 1; $r -le $t.Rows.Count; $r++) {
        for ($c = 1; $c -le $t.Columns.Count; $c++) {
            $cell = $t.Cell($r, $c)
            if ($cell.Range.Paragraphs.Count -ne 1) {
                continue
            }

            $para = $cell.Range.Paragraphs.Item(1)
            $styleName = $para.Range.Style.NameLocal

            # Cell range text always carries a trailing cell-mark
            # (CR + BEL); strip those control characters so we can
            # tell whether the cell is genuinely empty of content.
            $plainText = $para.Range.Text -replace "[\r\a]", ""

            if ($styleName -eq "Compact" -and $plainText.Length -eq 0) {
                $para.Range.Style = "Normal"
            }
        }
    }
}
